# Auto-generated edit script applying scheduled price/profit refresh to Asura_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 133: Big Brush, Big Dreams / Ginseng Angle Brush
$ws.Range("H133").Value = 52582.715
$ws.Range("J133").Value = 52582.715
$ws.Range("L133").Value = 52582.715
$ws.Range("N133").Value = -62702.715

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 13969.459
$ws.Range("I32").Value = 14140.418
$ws.Range("K32").Value = 14140.418
$ws.Range("M32").Value = -13853.418

# Row 63: Rivets Run through It / Mythrite Rivets
$ws.Range("H63").Value = 5803.625
$ws.Range("I63").Value = 5204.143
$ws.Range("J63").Value = 10000
$ws.Range("K63").Value = 5204.143
$ws.Range("L63").Value = 10000
$ws.Range("M63").Value = -4518.143
$ws.Range("N63").Value = -11372

# Row 66: A Riveting Revival (L) / Mythrite Rivets
$ws.Range("H66").Value = 5803.625
$ws.Range("I66").Value = 5204.143
$ws.Range("J66").Value = 10000
$ws.Range("K66").Value = 26020.715
$ws.Range("L66").Value = 50000
$ws.Range("M66").Value = -22588.715
$ws.Range("N66").Value = -56864

# Row 74: As the Bolt Flies / Titanium Nugget
$ws.Range("H74").Value = 1167.9445
$ws.Range("I74").Value = 931.375
$ws.Range("J74").Value = 1641.0834
$ws.Range("K74").Value = 931.375
$ws.Range("L74").Value = 1641.0834
$ws.Range("M74").Value = -57.375
$ws.Range("N74").Value = -3389.0834

# Row 77: Heavy Metal Banned (L) / Titanium Nugget
$ws.Range("H77").Value = 1167.9445
$ws.Range("I77").Value = 931.375
$ws.Range("J77").Value = 1641.0834
$ws.Range("K77").Value = 4656.875
$ws.Range("L77").Value = 8205.416999999999
$ws.Range("M77").Value = -288.875
$ws.Range("N77").Value = -16941.417

# Row 138: Don't Ask about the Rivets / Titanium Gold Helm of Casting
$ws.Range("H138").Value = 60773.77
$ws.Range("J138").Value = 60773.77
$ws.Range("L138").Value = 60773.77
$ws.Range("N138").Value = -71053.76999999999

# Row 141: Essays on Equipment / Ra'Kaznar Greaves of Maiming
$ws.Range("H141").Value = 55255.445
$ws.Range("J141").Value = 55255.445
$ws.Range("L141").Value = 55255.445
$ws.Range("N141").Value = -65615.44500000001

$ws = $wb.Worksheets.Item("BSM")
# Row 99: Meddle in Metal / Oroshigane Ingot
$ws.Range("H99").Value = 2835.1667
$ws.Range("I99").Value = 3000
$ws.Range("J99").Value = 2802.2
$ws.Range("K99").Value = 3000
$ws.Range("L99").Value = 2802.2
$ws.Range("M99").Value = -1502
$ws.Range("N99").Value = -5798.2

# Row 107: The Gold Experience / Deepgold Nugget
$ws.Range("H107").Value = 3961.1
$ws.Range("I107").Value = 4076.375
$ws.Range("K107").Value = 4076.375
$ws.Range("M107").Value = -2156.375

# Row 137: Dagger Swagger / Cobalt Tungsten Khukuri
$ws.Range("H137").Value = 63925.453
$ws.Range("J137").Value = 63925.453
$ws.Range("L137").Value = 63925.453
$ws.Range("N137").Value = -74125.45300000001

$ws = $wb.Worksheets.Item("CRP")
# Row 6: Got Your Back / Square Maple Shield
$ws.Range("H6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").Value = $null

# Row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 25003076
$ws.Range("I31").Value = 40002212
$ws.Range("J31").Value = 4513.3335
$ws.Range("K31").Value = 40002212
$ws.Range("L31").Value = 4513.3335
$ws.Range("M31").Value = -40001917
$ws.Range("N31").Value = -5103.3335

# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 25003076
$ws.Range("I34").Value = 40002212
$ws.Range("J34").Value = 4513.3335
$ws.Range("K34").Value = 40002212
$ws.Range("L34").Value = 4513.3335
$ws.Range("M34").Value = -40002010
$ws.Range("N34").Value = -4917.3335

# Row 41: The Lone Bowman / Oak Longbow
$ws.Range("H41").Value = 1750
$ws.Range("I41").Value = 1750
$ws.Range("K41").Value = 1750
$ws.Range("M41").Value = -1322

# Row 50: The Arsenal of Theocracy / Cobalt Halberd
$ws.Range("H50").Value = 35333.332
$ws.Range("J50").Value = 50500
$ws.Range("L50").Value = 50500
$ws.Range("N50").Value = -51750

# Row 51: Greenstone for Greenhorns / Jade Crook
$ws.Range("H51").Value = 34000
$ws.Range("J51").Value = 34000
$ws.Range("L51").Value = 34000
$ws.Range("N51").Value = -35472

# Row 60: Bowing to Greater Power / Yew Longbow
$ws.Range("H60").Value = 8930
$ws.Range("I60").Value = 8930
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 8930
$ws.Range("L60").Value = 0
$ws.Range("M60").Value = -8419
$ws.Range("N60").Value = $null

# Row 61: Incant Now, Think Later / Jade Crook
$ws.Range("H61").Value = 34000
$ws.Range("J61").Value = 34000
$ws.Range("L61").Value = 34000
$ws.Range("N61").Value = -34696

# Row 134: Wood You Be Quiet / Ceiba Lumber
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = $null
$ws.Range("N134").Value = $null

# Row 135: The Wing's Wings / Ceiba Wings
$ws.Range("H135").Value = 93014.44500000001
$ws.Range("J135").Value = 93014.44500000001
$ws.Range("L135").Value = 93014.44500000001
$ws.Range("N135").Value = -103154.445

# Row 137: Lament of the Lazylump / Dark Mahogany Fishing Rod
$ws.Range("H137").Value = 45995
$ws.Range("J137").Value = 61990
$ws.Range("L137").Value = 61990
$ws.Range("N137").Value = -72190

# Row 138: Bow Out / Acacia Longbow
$ws.Range("H138").Value = 43113.332
$ws.Range("J138").Value = 43113.332
$ws.Range("L138").Value = 43113.332
$ws.Range("N138").Value = -53393.332

$ws = $wb.Worksheets.Item("CUL")
# Row 5: What a Sap / Maple Syrup
$ws.Range("H5").Value = 1439.4138
$ws.Range("I5").Value = 2131.6428
$ws.Range("J5").Value = 793.3333
$ws.Range("K5").Value = 6394.928400000001
$ws.Range("L5").Value = 2379.9999
$ws.Range("M5").Value = -6282.928400000001
$ws.Range("N5").Value = -2603.9999

# Row 12: Butter Me Up / Kukuru Butter
$ws.Range("H12").Value = 920229.7
$ws.Range("I12").Value = 109
$ws.Range("J12").Value = 1136728.6
$ws.Range("K12").Value = 327
$ws.Range("L12").Value = 3410185.8
$ws.Range("M12").Value = -154
$ws.Range("N12").Value = -3410531.8

# Row 131: The Mountain Steeped / Tsai tou Vounou
$ws.Range("H131").Value = 14708723
$ws.Range("I131").Value = 14671.429
$ws.Range("J131").Value = 16394926
$ws.Range("K131").Value = 44014.287
$ws.Range("L131").Value = 49184778
$ws.Range("M131").Value = -38974.287
$ws.Range("N131").Value = -49194858

# Row 132: More Mezcal / Cooking Mezcal
$ws.Range("H132").Value = 1544.2174
$ws.Range("I132").Value = 950.625
$ws.Range("J132").Value = 1860.8
$ws.Range("K132").Value = 8555.625
$ws.Range("L132").Value = 16747.2
$ws.Range("M132").Value = -6025.625
$ws.Range("N132").Value = -21807.2

# Row 135: Not-so-secret Ingredient / Royal Maple Syrup
$ws.Range("H135").Value = 1439.4138
$ws.Range("I135").Value = 2131.6428
$ws.Range("J135").Value = 793.3333
$ws.Range("K135").Value = 19184.7852
$ws.Range("L135").Value = 7139.9997
$ws.Range("M135").Value = -16649.7852
$ws.Range("N135").Value = -12209.9997

$ws = $wb.Worksheets.Item("LTW")
# Row 134: Freezing Fingers / Crocodileskin Fingerless Gloves of Striking
$ws.Range("H134").Value = 81214.5
$ws.Range("J134").Value = 81214.5
$ws.Range("L134").Value = 81214.5
$ws.Range("N134").Value = -91354.5

# Row 140: Worqor Zormor or Bust / Gargantuaskin Shoes of Healing
$ws.Range("H140").Value = 75262
$ws.Range("J140").Value = 75262
$ws.Range("L140").Value = 75262
$ws.Range("N140").Value = -85622

$ws = $wb.Worksheets.Item("WVR")
# Row 110: Suits You / Iridescent Acton of Aiming
$ws.Range("H110").Value = 30644
$ws.Range("J110").Value = 30644
$ws.Range("L110").Value = 30644
$ws.Range("N110").Value = -38824

# Row 138: Halfgloves, Full Effort / Rroneek Serge Halfgloves of Healing
$ws.Range("H138").Value = 84625
$ws.Range("J138").Value = 84625
$ws.Range("L138").Value = 84625
$ws.Range("N138").Value = -94905

# Row 140: Glamorous Gloves / Thunderyards Silk Gloves of Casting
$ws.Range("H140").Value = 43166.5
$ws.Range("J140").Value = 43166.5
$ws.Range("L140").Value = 43166.5
$ws.Range("N140").Value = -53526.5
